$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("status_label") - shifts NCTId..results from B..I to C..J
$ws.Columns("B").Insert()

# Set new header for status_label
$ws.Range("B1").Value = "status_label"

# Fill status_label (column B) for every data row, derived from the status (column A)
$ws.Range("B2").Value = "rouge"
$ws.Range("B3").Value = "rouge"
$ws.Range("B4").Value = "rouge"
$ws.Range("B5").Value = "rouge"
$ws.Range("B6").Value = "rouge"
$ws.Range("B7").Value = "rouge"
$ws.Range("B8").Value = "rouge"
$ws.Range("B9").Value = "orange"
$ws.Range("B10").Value = "rouge"
$ws.Range("B11").Value = "rouge"
$ws.Range("B12").Value = "rouge"
$ws.Range("B13").Value = "rouge"
$ws.Range("B14").Value = "rouge"
$ws.Range("B15").Value = "rouge"
$ws.Range("B16").Value = "rouge"
$ws.Range("B17").Value = "rouge"
$ws.Range("B18").Value = "rouge"
$ws.Range("B19").Value = "rouge"
$ws.Range("B20").Value = "rouge"
$ws.Range("B21").Value = "vert"
$ws.Range("B22").Value = "rouge"
$ws.Range("B23").Value = "rouge"
$ws.Range("B24").Value = "vert"
$ws.Range("B25").Value = "rouge"
$ws.Range("B26").Value = "rouge"
$ws.Range("B27").Value = "rouge"

# Re-sync rows whose record order changed (status/NCTId/title/acronym/results realigned)
$ws.Range("A7").Value = "🟥"
$ws.Range("C7").Value = "NCT00870909"
$ws.Range("F7").Value = "Anodal & Cathodal tDCS for Treatment of Resistant Auditory Hallucinations in Schizophrenia"
$ws.Range("G7").Value = ""
$ws.Range("H7").Value = $false
$ws.Range("I7").Value = $false
$ws.Range("J7").Value = $true

$ws.Range("A8").Value = "🟥"
$ws.Range("C8").Value = "NCT02717260"
$ws.Range("F8").Value = "Inhibition Control Modulation by Transcranial Random Noise Stimulation (tRNS ) on the Prefrontal Cortex Measured by Change in Go no Test"
$ws.Range("G8").Value = "inhibistim"
$ws.Range("H8").Value = $false
$ws.Range("I8").Value = $false
$ws.Range("J8").Value = $false

$ws.Range("A9").Value = "🟧"
$ws.Range("C9").Value = "NCT02438163"
$ws.Range("F9").Value = "Study of Neuroplasticity on Depressed Patients Versus Healthy Subjects : Modulation of the MEP Size Induced by Theta Burst Stimulation"
$ws.Range("G9").Value = "DEPLAS"
$ws.Range("H9").Value = $false
$ws.Range("I9").Value = $true
$ws.Range("J9").Value = $true

$ws.Range("A10").Value = "🟥"
$ws.Range("C10").Value = "NCT01891929"
$ws.Range("F10").Value = "Specific Cognitive Remediation for Schizophrenia (RECOS) and Sheltered Employment: a Multicentre Controlled Randomized Trial."
$ws.Range("G10").Value = "RemedRehab"
$ws.Range("H10").Value = $false
$ws.Range("I10").Value = $false
$ws.Range("J10").Value = $false

$ws.Range("A11").Value = "🟥"
$ws.Range("C11").Value = "NCT02734927"
$ws.Range("F11").Value = "Motivation and Executive Control in Schizophrenia"
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = $false
$ws.Range("I11").Value = $false
$ws.Range("J11").Value = $false

$ws.Range("A12").Value = "🟥"
$ws.Range("C12").Value = "NCT02667834"
$ws.Range("F12").Value = "Efficiency of the French Translation Program:Social Cognition and Interactive Training (SCIT) of Negative Symptoms in Schizophrenia"
$ws.Range("G12").Value = "SCIT-VF"
$ws.Range("H12").Value = $false
$ws.Range("I12").Value = $false
$ws.Range("J12").Value = $false

$ws.Range("A13").Value = "🟥"
$ws.Range("C13").Value = "NCT02949453"
$ws.Range("F13").Value = "ACCeptation and Qualitative Evaluation of Phone-delivered Intervention To Prevent Suicide Reattempt"
$ws.Range("G13").Value = "ACCEPT-S"
$ws.Range("H13").Value = $false
$ws.Range("I13").Value = $false
$ws.Range("J13").Value = $false

$ws.Range("A14").Value = "🟥"
$ws.Range("C14").Value = "NCT03492970"
$ws.Range("F14").Value = "Characterization of Behavioral Disorders and 24 H-melatonin Level in Adults With Smith Magenis Syndrome"
$ws.Range("G14").Value = "SMS-adults"
$ws.Range("H14").Value = $false
$ws.Range("I14").Value = $false
$ws.Range("J14").Value = $false

$ws.Range("A15").Value = "🟥"
$ws.Range("C15").Value = "NCT03028545"
$ws.Range("F15").Value = "Representations, Strategies and Identity Redefinition in the Recovery Process: Exploratory Study"
$ws.Range("G15").Value = "EPR"
$ws.Range("H15").Value = $false
$ws.Range("I15").Value = $false
$ws.Range("J15").Value = $false

$ws.Range("A16").Value = "🟥"
$ws.Range("C16").Value = "NCT03688516"
$ws.Range("F16").Value = "Effects of Emotion on Episodic Memory in Typically Developing Children and Children With Williams-Beuren Syndrome"
$ws.Range("G16").Value = "EEM-TAdev"
$ws.Range("H16").Value = $false
$ws.Range("I16").Value = $false
$ws.Range("J16").Value = $false

$ws.Range("A17").Value = "🟥"
$ws.Range("C17").Value = "NCT04141540"
$ws.Range("F17").Value = "Translational 22q11.2:`"Molecular Variants Associated With Schizophrenia: Differential Analysis of Monozygotic Twins With Variable Phenotypic 22q11.2 Microdeletional Syndrom`""
$ws.Range("G17").Value = "CSRK05"
$ws.Range("H17").Value = $false
$ws.Range("I17").Value = $false
$ws.Range("J17").Value = $false

$ws.Range("A20").Value = "🟥"
$ws.Range("C20").Value = "NCT04147988"
$ws.Range("F20").Value = "Describing, Detecting and Orienting Adults Without Intellectual Disability Asking for an Autism Spectrum Disorder Diagnosis"
$ws.Range("G20").Value = "ASSORT"
$ws.Range("H20").Value = $false
$ws.Range("I20").Value = $false
$ws.Range("J20").Value = $false

$ws.Range("A21").Value = "🟩"
$ws.Range("C21").Value = "NCT02608333"
$ws.Range("F21").Value = "Impact of Early Intervention on the Global Development of Children With Autism Spectrum Disorder in a European French-speaking Population Dr Marie-Maude GEOFFRAY Le Vinatier Hospital"
$ws.Range("G21").Value = "IDEA"
$ws.Range("H21").Value = $true
$ws.Range("I21").Value = $true
$ws.Range("J21").Value = $true

$ws.Range("A25").Value = "🟥"
$ws.Range("C25").Value = "NCT05213442"
$ws.Range("F25").Value = "Liver Status in Psychiatric Inpatients"
$ws.Range("G25").Value = "LIVERSPIN"
$ws.Range("H25").Value = $false
$ws.Range("I25").Value = $false
$ws.Range("J25").Value = $false

$ws.Range("A26").Value = "🟥"
$ws.Range("C26").Value = "NCT04245306"
$ws.Range("F26").Value = "Pilot Evaluation of a New Computerized Test for Pragmatic Inferences in Children with Autism Spectrum Disorder Aged 8-12 Years Old"
$ws.Range("G26").Value = "TIPI"
$ws.Range("H26").Value = $false
$ws.Range("I26").Value = $false
$ws.Range("J26").Value = $false
